$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row after the header row; Excel automatically copies the
# formatting of the row above (the bold header row) onto the new row.
$ws.Rows.Item(2).Insert()

# Reset the row-4 (previously row 5) custom height that belonged to the
# old "Code challenges" row back to the sheet default, since that content
# moved/changed.
$ws.Rows.Item(6).AutoFit()

# ---- Row 2 (new) : 8:00 - 8:15 / Outline agenda, introduction / Patrick or Cory
$ws.Cells.Item(2,1).Value = "8:00 - 8:15"
$ws.Cells.Item(2,2).Value = "Outline agenda, introduction"
$ws.Cells.Item(2,3).Value = "Patrick or Cory"
# This row was seeded from the bold header formatting; un-bold the time
# and presenter cells but keep the normal wrapped body style for column B.
$ws.Cells.Item(2,1).Font.Bold = $false
$ws.Cells.Item(2,2).Font.Bold = $false
$ws.Cells.Item(2,2).WrapText = $true
$ws.Cells.Item(2,2).VerticalAlignment = -4108
$ws.Cells.Item(2,3).Font.Bold = $false

# ---- Row 3 : 8:15 - 8:45 / Presentation: Intro to refactoring, How, When / Ryan
$ws.Cells.Item(3,1).Value = "8:15 - 8:45"
$ws.Cells.Item(3,2).Value = "Presentation: Intro to refactoring, How, When"
$ws.Cells.Item(3,3).Value = "Ryan"

# ---- Row 4 : 8:45 - 9:45 / Code challenges and discussion - ... / All
$ws.Cells.Item(4,1).Value = "8:45 - 9:45"
$ws.Cells.Item(4,2).Value = "Code challenges and discussion - Review piece of code as a group and identify bad practices. Use this as a lead-in to present and discuss specific anti-patterns and code smells"
$ws.Cells.Item(4,3).Value = "All"
$ws.Rows.Item(4).RowHeight = 60

# ---- Row 5 : 9:45 - 10:00 / Break
$ws.Cells.Item(5,1).Value = "9:45 - 10:00"
$ws.Cells.Item(5,2).Value = "Break"
$ws.Cells.Item(5,3).ClearContents()

# ---- Row 6 : 10:00 - 10:45 / Refactoring related portions of Clean Code talk / Cory
$ws.Cells.Item(6,1).Value = "10:00 - 10:45"
$ws.Cells.Item(6,2).Value = "Refactoring related portions of Clean Code talk"
$ws.Cells.Item(6,3).Value = "Cory"

# ---- Row 7 : 10:45 - 12:00 / Continue challenges and discussions / All
$ws.Cells.Item(7,1).Value = "10:45 - 12:00"
$ws.Cells.Item(7,2).Value = "Continue challenges and discussions"
$ws.Cells.Item(7,3).Value = "All"

# ---- Row 8 : 12:00 - 01:00 / Lunch
$ws.Cells.Item(8,1).Value = "12:00 - 01:00"
$ws.Cells.Item(8,2).Value = "Lunch"
$ws.Cells.Item(8,3).ClearContents()

# ---- Row 9 : 01:00 - 01:30 / Presentation: ?Specific refactoring techniques? / Patrick?
$ws.Cells.Item(9,1).Value = "01:00 - 01:30"
$ws.Cells.Item(9,2).Value = "Presentation: ?Specific refactoring techniques?"
$ws.Cells.Item(9,3).Value = "Patrick?"

# ---- Row 10 : 01:30 - 03:00 / 1 on 1 code review and refactoring assistance / All
$ws.Cells.Item(10,1).Value = "01:30 - 03:00"
$ws.Cells.Item(10,2).Value = "1 on 1 code review and refactoring assistance"
$ws.Cells.Item(10,3).Value = "All"

# ---- Row 11 : 03:00 - 04:00 / Review samples from 1 on 1 code reviews with the group. / All?
$ws.Cells.Item(11,1).Value = "03:00 - 04:00"
$ws.Cells.Item(11,2).Value = "Review samples from 1 on 1 code reviews with the group."
$ws.Cells.Item(11,3).Value = "All?"
$ws.Rows.Item(11).RowHeight = 30

# ---- Row 12 (new) : 4:00 - 5:00 / Review refactored example application
$ws.Cells.Item(12,1).Value = "4:00 - 5:00"
$ws.Cells.Item(12,2).Value = "Review refactored example application"
$ws.Cells.Item(12,3).ClearContents()

$ws.Range("B13").Select() | Out-Null
